$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# New GPS-tracking entries for two extra training days (J-2, J-1) and one
# new match ("N3 J6 VS Riviera"), appended as sheet rows 698-717.
# ----------------------------------------------------------------------

# Step 1: "Temps joue" (col G) values, written row by row, first -- this
# matches the order new shared strings were introduced in the saved file.
$ws.Cells.Item(698, 7).Value = "01:12:43"
$ws.Cells.Item(699, 7).Value = "01:11:52"
$ws.Cells.Item(700, 7).Value = "01:12:58"
$ws.Cells.Item(701, 7).Value = "01:12:35"
$ws.Cells.Item(702, 7).Value = "00:18:07"
$ws.Cells.Item(703, 7).Value = "00:19:26"
$ws.Cells.Item(704, 7).Value = "00:19:10"
$ws.Cells.Item(705, 7).Value = "01:42:25"
$ws.Cells.Item(706, 7).Value = "01:41:55"
$ws.Cells.Item(707, 7).Value = "01:26:19"
$ws.Cells.Item(708, 7).Value = "01:42:47"
$ws.Cells.Item(709, 7).Value = "01:42:47"
$ws.Cells.Item(710, 7).Value = "00:15:51"
$ws.Cells.Item(711, 7).Value = "01:42:47"
$ws.Cells.Item(712, 7).Value = "00:07:48"
$ws.Cells.Item(713, 7).Value = "01:42:47"
$ws.Cells.Item(714, 7).Value = "01:05:38"
$ws.Cells.Item(715, 7).Value = "01:34:37"
$ws.Cells.Item(716, 7).Value = "01:41:55"
$ws.Cells.Item(717, 7).Value = "00:35:25"

# Step 2: "Type" (col A) for the new match rows 705-717 -- written after all
# the "Temps joue" values so its shared string is appended last, like in the
# saved workbook.
$ws.Cells.Item(705, 1).Value = "N3 J6 VS Riviera"
$ws.Cells.Item(706, 1).Value = "N3 J6 VS Riviera"
$ws.Cells.Item(707, 1).Value = "N3 J6 VS Riviera"
$ws.Cells.Item(708, 1).Value = "N3 J6 VS Riviera"
$ws.Cells.Item(709, 1).Value = "N3 J6 VS Riviera"
$ws.Cells.Item(710, 1).Value = "N3 J6 VS Riviera"
$ws.Cells.Item(711, 1).Value = "N3 J6 VS Riviera"
$ws.Cells.Item(712, 1).Value = "N3 J6 VS Riviera"
$ws.Cells.Item(713, 1).Value = "N3 J6 VS Riviera"
$ws.Cells.Item(714, 1).Value = "N3 J6 VS Riviera"
$ws.Cells.Item(715, 1).Value = "N3 J6 VS Riviera"
$ws.Cells.Item(716, 1).Value = "N3 J6 VS Riviera"
$ws.Cells.Item(717, 1).Value = "N3 J6 VS Riviera"

# Step 3: "Type" (col A) for the training rows 698-704 (reuses the existing
# "Entrainement" string already in the workbook).
$ws.Cells.Item(698, 1).Value = "Entrainement"
$ws.Cells.Item(699, 1).Value = "Entrainement"
$ws.Cells.Item(700, 1).Value = "Entrainement"
$ws.Cells.Item(701, 1).Value = "Entrainement"
$ws.Cells.Item(702, 1).Value = "Entrainement"
$ws.Cells.Item(703, 1).Value = "Entrainement"
$ws.Cells.Item(704, 1).Value = "Entrainement"

# Step 4: remaining columns B-F (Date/Periode/MD/Nom du joueur/Poste) and
# H-V (the numeric GPS metrics).
$ws.Cells.Item(698, 2).Value = 45946.0
$ws.Cells.Item(698, 3).Value = "Global"
$ws.Cells.Item(698, 4).Value = "J-2"
$ws.Cells.Item(698, 5).Value = "Mattheo Haon"
$ws.Cells.Item(698, 6).Value = "right back"
$ws.Cells.Item(698, 8).Value = 3.37
$ws.Cells.Item(698, 9).Value = 0.12
$ws.Cells.Item(698, 10).Value = 3.25
$ws.Cells.Item(698, 11).Value = 0.09
$ws.Cells.Item(698, 12).Value = 0.03
$ws.Cells.Item(698, 13).Value = 0.0
$ws.Cells.Item(698, 14).Value = 0.0
$ws.Cells.Item(698, 15).Value = 0.0
$ws.Cells.Item(698, 16).Value = 2.69
$ws.Cells.Item(698, 17).Value = 23.34
$ws.Cells.Item(698, 18).Value = 5.2
$ws.Cells.Item(698, 19).Value = 9.0
$ws.Cells.Item(698, 20).Value = 8.0
$ws.Cells.Item(698, 21).Value = 3.0
$ws.Cells.Item(698, 22).Value = 0.0

$ws.Cells.Item(699, 2).Value = 45946.0
$ws.Cells.Item(699, 3).Value = "Global"
$ws.Cells.Item(699, 4).Value = "J-2"
$ws.Cells.Item(699, 5).Value = "Malik Boussaid"
$ws.Cells.Item(699, 6).Value = "right back"
$ws.Cells.Item(699, 8).Value = 3.28
$ws.Cells.Item(699, 9).Value = 0.04
$ws.Cells.Item(699, 10).Value = 3.24
$ws.Cells.Item(699, 11).Value = 0.04
$ws.Cells.Item(699, 12).Value = 0.0
$ws.Cells.Item(699, 13).Value = 0.0
$ws.Cells.Item(699, 14).Value = 0.0
$ws.Cells.Item(699, 15).Value = 0.0
$ws.Cells.Item(699, 16).Value = 2.15
$ws.Cells.Item(699, 17).Value = 20.43
$ws.Cells.Item(699, 18).Value = 4.54
$ws.Cells.Item(699, 19).Value = 19.0
$ws.Cells.Item(699, 20).Value = 1.0
$ws.Cells.Item(699, 21).Value = 19.0
$ws.Cells.Item(699, 22).Value = 5.0

$ws.Cells.Item(700, 2).Value = 45946.0
$ws.Cells.Item(700, 3).Value = "Global"
$ws.Cells.Item(700, 4).Value = "J-2"
$ws.Cells.Item(700, 5).Value = "Emmanuel Valey"
$ws.Cells.Item(700, 6).Value = "left forward"
$ws.Cells.Item(700, 8).Value = 3.31
$ws.Cells.Item(700, 9).Value = 0.01
$ws.Cells.Item(700, 10).Value = 3.3
$ws.Cells.Item(700, 11).Value = 0.01
$ws.Cells.Item(700, 12).Value = 0.0
$ws.Cells.Item(700, 13).Value = 0.0
$ws.Cells.Item(700, 14).Value = 0.0
$ws.Cells.Item(700, 15).Value = 0.0
$ws.Cells.Item(700, 16).Value = 1.94
$ws.Cells.Item(700, 17).Value = 17.11
$ws.Cells.Item(700, 18).Value = 3.54
$ws.Cells.Item(700, 19).Value = 5.0
$ws.Cells.Item(700, 20).Value = 0.0
$ws.Cells.Item(700, 21).Value = 2.0
$ws.Cells.Item(700, 22).Value = 0.0

$ws.Cells.Item(701, 2).Value = 45946.0
$ws.Cells.Item(701, 3).Value = "Global"
$ws.Cells.Item(701, 4).Value = "J-2"
$ws.Cells.Item(701, 5).Value = "Karahali Souaré"
$ws.Cells.Item(701, 6).Value = "right forward"
$ws.Cells.Item(701, 8).Value = 2.86
$ws.Cells.Item(701, 9).Value = 0.02
$ws.Cells.Item(701, 10).Value = 2.84
$ws.Cells.Item(701, 11).Value = 0.02
$ws.Cells.Item(701, 12).Value = 0.0
$ws.Cells.Item(701, 13).Value = 0.0
$ws.Cells.Item(701, 14).Value = 0.0
$ws.Cells.Item(701, 15).Value = 0.0
$ws.Cells.Item(701, 16).Value = 2.0
$ws.Cells.Item(701, 17).Value = 18.66
$ws.Cells.Item(701, 18).Value = 4.3
$ws.Cells.Item(701, 19).Value = 17.0
$ws.Cells.Item(701, 20).Value = 1.0
$ws.Cells.Item(701, 21).Value = 17.0
$ws.Cells.Item(701, 22).Value = 1.0

$ws.Cells.Item(702, 2).Value = 45947.0
$ws.Cells.Item(702, 3).Value = "Global"
$ws.Cells.Item(702, 4).Value = "J-1"
$ws.Cells.Item(702, 5).Value = "Karahali Souaré"
$ws.Cells.Item(702, 6).Value = "right forward"
$ws.Cells.Item(702, 8).Value = 1.68
$ws.Cells.Item(702, 9).Value = 0.13
$ws.Cells.Item(702, 10).Value = 1.54
$ws.Cells.Item(702, 11).Value = 0.11
$ws.Cells.Item(702, 12).Value = 0.01
$ws.Cells.Item(702, 13).Value = 0.01
$ws.Cells.Item(702, 14).Value = 0.0
$ws.Cells.Item(702, 15).Value = 1.0
$ws.Cells.Item(702, 16).Value = 5.32
$ws.Cells.Item(702, 17).Value = 30.12
$ws.Cells.Item(702, 18).Value = 5.13
$ws.Cells.Item(702, 19).Value = 29.0
$ws.Cells.Item(702, 20).Value = 4.0
$ws.Cells.Item(702, 21).Value = 18.0
$ws.Cells.Item(702, 22).Value = 10.0

$ws.Cells.Item(703, 2).Value = 45947.0
$ws.Cells.Item(703, 3).Value = "Global"
$ws.Cells.Item(703, 4).Value = "J-1"
$ws.Cells.Item(703, 5).Value = "Malik Boussaid"
$ws.Cells.Item(703, 6).Value = "right back"
$ws.Cells.Item(703, 8).Value = 2.07
$ws.Cells.Item(703, 9).Value = 0.2
$ws.Cells.Item(703, 10).Value = 1.86
$ws.Cells.Item(703, 11).Value = 0.15
$ws.Cells.Item(703, 12).Value = 0.05
$ws.Cells.Item(703, 13).Value = 0.0
$ws.Cells.Item(703, 14).Value = 0.0
$ws.Cells.Item(703, 15).Value = 0.0
$ws.Cells.Item(703, 16).Value = 6.05
$ws.Cells.Item(703, 17).Value = 24.08
$ws.Cells.Item(703, 18).Value = 4.86
$ws.Cells.Item(703, 19).Value = 19.0
$ws.Cells.Item(703, 20).Value = 4.0
$ws.Cells.Item(703, 21).Value = 22.0
$ws.Cells.Item(703, 22).Value = 2.0

$ws.Cells.Item(704, 2).Value = 45947.0
$ws.Cells.Item(704, 3).Value = "Global"
$ws.Cells.Item(704, 4).Value = "J-1"
$ws.Cells.Item(704, 5).Value = "Mattheo Haon"
$ws.Cells.Item(704, 6).Value = "right back"
$ws.Cells.Item(704, 8).Value = 1.89
$ws.Cells.Item(704, 9).Value = 0.19
$ws.Cells.Item(704, 10).Value = 1.69
$ws.Cells.Item(704, 11).Value = 0.14
$ws.Cells.Item(704, 12).Value = 0.05
$ws.Cells.Item(704, 13).Value = 0.0
$ws.Cells.Item(704, 14).Value = 0.0
$ws.Cells.Item(704, 15).Value = 1.0
$ws.Cells.Item(704, 16).Value = 5.91
$ws.Cells.Item(704, 17).Value = 25.52
$ws.Cells.Item(704, 18).Value = 4.4
$ws.Cells.Item(704, 19).Value = 9.0
$ws.Cells.Item(704, 20).Value = 2.0
$ws.Cells.Item(704, 21).Value = 13.0
$ws.Cells.Item(704, 22).Value = 1.0

$ws.Cells.Item(705, 2).Value = 45948.0
$ws.Cells.Item(705, 3).Value = "Global"
$ws.Cells.Item(705, 4).Value = "M"
$ws.Cells.Item(705, 5).Value = "Romain Thunet"
$ws.Cells.Item(705, 6).Value = "center back"
$ws.Cells.Item(705, 8).Value = 10.42
$ws.Cells.Item(705, 9).Value = 1.27
$ws.Cells.Item(705, 10).Value = 9.13
$ws.Cells.Item(705, 11).Value = 0.93
$ws.Cells.Item(705, 12).Value = 0.31
$ws.Cells.Item(705, 13).Value = 0.05
$ws.Cells.Item(705, 14).Value = 0.0
$ws.Cells.Item(705, 15).Value = 5.0
$ws.Cells.Item(705, 16).Value = 6.05
$ws.Cells.Item(705, 17).Value = 27.36
$ws.Cells.Item(705, 18).Value = 4.3
$ws.Cells.Item(705, 19).Value = 32.0
$ws.Cells.Item(705, 20).Value = 3.0
$ws.Cells.Item(705, 21).Value = 28.0
$ws.Cells.Item(705, 22).Value = 10.0

$ws.Cells.Item(706, 2).Value = 45948.0
$ws.Cells.Item(706, 3).Value = "Global"
$ws.Cells.Item(706, 4).Value = "M"
$ws.Cells.Item(706, 5).Value = "Yoann Martelat"
$ws.Cells.Item(706, 6).Value = "center midfield"
$ws.Cells.Item(706, 8).Value = 12.42
$ws.Cells.Item(706, 9).Value = 2.87
$ws.Cells.Item(706, 10).Value = 9.51
$ws.Cells.Item(706, 11).Value = 2.31
$ws.Cells.Item(706, 12).Value = 0.56
$ws.Cells.Item(706, 13).Value = 0.04
$ws.Cells.Item(706, 14).Value = 0.0
$ws.Cells.Item(706, 15).Value = 3.0
$ws.Cells.Item(706, 16).Value = 7.28
$ws.Cells.Item(706, 17).Value = 27.85
$ws.Cells.Item(706, 18).Value = 5.09
$ws.Cells.Item(706, 19).Value = 32.0
$ws.Cells.Item(706, 20).Value = 2.0
$ws.Cells.Item(706, 21).Value = 43.0
$ws.Cells.Item(706, 22).Value = 7.0

$ws.Cells.Item(707, 2).Value = 45948.0
$ws.Cells.Item(707, 3).Value = "Global"
$ws.Cells.Item(707, 4).Value = "M"
$ws.Cells.Item(707, 5).Value = "Sofiane Belle"
$ws.Cells.Item(707, 6).Value = "left forward"
$ws.Cells.Item(707, 8).Value = 8.78
$ws.Cells.Item(707, 9).Value = 1.61
$ws.Cells.Item(707, 10).Value = 7.15
$ws.Cells.Item(707, 11).Value = 0.93
$ws.Cells.Item(707, 12).Value = 0.45
$ws.Cells.Item(707, 13).Value = 0.25
$ws.Cells.Item(707, 14).Value = 0.0
$ws.Cells.Item(707, 15).Value = 19.0
$ws.Cells.Item(707, 16).Value = 5.98
$ws.Cells.Item(707, 17).Value = 29.81
$ws.Cells.Item(707, 18).Value = 4.25
$ws.Cells.Item(707, 19).Value = 32.0
$ws.Cells.Item(707, 20).Value = 3.0
$ws.Cells.Item(707, 21).Value = 27.0
$ws.Cells.Item(707, 22).Value = 8.0

$ws.Cells.Item(708, 2).Value = 45948.0
$ws.Cells.Item(708, 3).Value = "Global"
$ws.Cells.Item(708, 4).Value = "M"
$ws.Cells.Item(708, 5).Value = "Malik Boussaid"
$ws.Cells.Item(708, 6).Value = "right back"
$ws.Cells.Item(708, 8).Value = 12.26
$ws.Cells.Item(708, 9).Value = 2.48
$ws.Cells.Item(708, 10).Value = 9.75
$ws.Cells.Item(708, 11).Value = 1.67
$ws.Cells.Item(708, 12).Value = 0.72
$ws.Cells.Item(708, 13).Value = 0.12
$ws.Cells.Item(708, 14).Value = 0.0
$ws.Cells.Item(708, 15).Value = 17.0
$ws.Cells.Item(708, 16).Value = 6.97
$ws.Cells.Item(708, 17).Value = 27.66
$ws.Cells.Item(708, 18).Value = 5.69
$ws.Cells.Item(708, 19).Value = 46.0
$ws.Cells.Item(708, 20).Value = 6.0
$ws.Cells.Item(708, 21).Value = 32.0
$ws.Cells.Item(708, 22).Value = 16.0

$ws.Cells.Item(709, 2).Value = 45948.0
$ws.Cells.Item(709, 3).Value = "Global"
$ws.Cells.Item(709, 4).Value = "M"
$ws.Cells.Item(709, 5).Value = "Mattheo Haon"
$ws.Cells.Item(709, 6).Value = "right back"
$ws.Cells.Item(709, 8).Value = 10.9
$ws.Cells.Item(709, 9).Value = 1.63
$ws.Cells.Item(709, 10).Value = 9.25
$ws.Cells.Item(709, 11).Value = 1.2
$ws.Cells.Item(709, 12).Value = 0.41
$ws.Cells.Item(709, 13).Value = 0.04
$ws.Cells.Item(709, 14).Value = 0.0
$ws.Cells.Item(709, 15).Value = 3.0
$ws.Cells.Item(709, 16).Value = 6.29
$ws.Cells.Item(709, 17).Value = 27.15
$ws.Cells.Item(709, 18).Value = 4.37
$ws.Cells.Item(709, 19).Value = 25.0
$ws.Cells.Item(709, 20).Value = 2.0
$ws.Cells.Item(709, 21).Value = 26.0
$ws.Cells.Item(709, 22).Value = 7.0

$ws.Cells.Item(710, 2).Value = 45948.0
$ws.Cells.Item(710, 3).Value = "Global"
$ws.Cells.Item(710, 4).Value = "M"
$ws.Cells.Item(710, 5).Value = "Ilyes Boughanmi"
$ws.Cells.Item(710, 6).Value = "center forward"
$ws.Cells.Item(710, 8).Value = 1.75
$ws.Cells.Item(710, 9).Value = 0.38
$ws.Cells.Item(710, 10).Value = 1.37
$ws.Cells.Item(710, 11).Value = 0.26
$ws.Cells.Item(710, 12).Value = 0.1
$ws.Cells.Item(710, 13).Value = 0.02
$ws.Cells.Item(710, 14).Value = 0.0
$ws.Cells.Item(710, 15).Value = 2.0
$ws.Cells.Item(710, 16).Value = 6.62
$ws.Cells.Item(710, 17).Value = 27.0
$ws.Cells.Item(710, 18).Value = 3.45
$ws.Cells.Item(710, 19).Value = 5.0
$ws.Cells.Item(710, 20).Value = 0.0
$ws.Cells.Item(710, 21).Value = 6.0
$ws.Cells.Item(710, 22).Value = 5.0

$ws.Cells.Item(711, 2).Value = 45948.0
$ws.Cells.Item(711, 3).Value = "Global"
$ws.Cells.Item(711, 4).Value = "M"
$ws.Cells.Item(711, 5).Value = "Naim Dhib"
$ws.Cells.Item(711, 6).Value = "center midfield"
$ws.Cells.Item(711, 8).Value = 10.75
$ws.Cells.Item(711, 9).Value = 1.55
$ws.Cells.Item(711, 10).Value = 9.19
$ws.Cells.Item(711, 11).Value = 1.09
$ws.Cells.Item(711, 12).Value = 0.43
$ws.Cells.Item(711, 13).Value = 0.06
$ws.Cells.Item(711, 14).Value = 0.0
$ws.Cells.Item(711, 15).Value = 9.0
$ws.Cells.Item(711, 16).Value = 6.17
$ws.Cells.Item(711, 17).Value = 28.51
$ws.Cells.Item(711, 18).Value = 4.73
$ws.Cells.Item(711, 19).Value = 47.0
$ws.Cells.Item(711, 20).Value = 2.0
$ws.Cells.Item(711, 21).Value = 30.0
$ws.Cells.Item(711, 22).Value = 12.0

$ws.Cells.Item(712, 2).Value = 45948.0
$ws.Cells.Item(712, 3).Value = "Global"
$ws.Cells.Item(712, 4).Value = "M"
$ws.Cells.Item(712, 5).Value = "Jeremie Laurent"
$ws.Cells.Item(712, 6).Value = "left forward"
$ws.Cells.Item(712, 8).Value = 1.08
$ws.Cells.Item(712, 9).Value = 0.32
$ws.Cells.Item(712, 10).Value = 0.76
$ws.Cells.Item(712, 11).Value = 0.24
$ws.Cells.Item(712, 12).Value = 0.07
$ws.Cells.Item(712, 13).Value = 0.02
$ws.Cells.Item(712, 14).Value = 0.0
$ws.Cells.Item(712, 15).Value = 1.0
$ws.Cells.Item(712, 16).Value = 8.32
$ws.Cells.Item(712, 17).Value = 27.71
$ws.Cells.Item(712, 18).Value = 4.39
$ws.Cells.Item(712, 19).Value = 11.0
$ws.Cells.Item(712, 20).Value = 3.0
$ws.Cells.Item(712, 21).Value = 4.0
$ws.Cells.Item(712, 22).Value = 1.0

$ws.Cells.Item(713, 2).Value = 45948.0
$ws.Cells.Item(713, 3).Value = "Global"
$ws.Cells.Item(713, 4).Value = "M"
$ws.Cells.Item(713, 5).Value = "Naim Ighbane"
$ws.Cells.Item(713, 6).Value = "center back"
$ws.Cells.Item(713, 8).Value = 10.3
$ws.Cells.Item(713, 9).Value = 1.41
$ws.Cells.Item(713, 10).Value = 8.87
$ws.Cells.Item(713, 11).Value = 0.95
$ws.Cells.Item(713, 12).Value = 0.37
$ws.Cells.Item(713, 13).Value = 0.11
$ws.Cells.Item(713, 14).Value = 0.0
$ws.Cells.Item(713, 15).Value = 8.0
$ws.Cells.Item(713, 16).Value = 5.92
$ws.Cells.Item(713, 17).Value = 28.79
$ws.Cells.Item(713, 18).Value = 4.92
$ws.Cells.Item(713, 19).Value = 31.0
$ws.Cells.Item(713, 20).Value = 5.0
$ws.Cells.Item(713, 21).Value = 34.0
$ws.Cells.Item(713, 22).Value = 8.0

$ws.Cells.Item(714, 2).Value = 45948.0
$ws.Cells.Item(714, 3).Value = "Global"
$ws.Cells.Item(714, 4).Value = "M"
$ws.Cells.Item(714, 5).Value = "Karahali Souaré"
$ws.Cells.Item(714, 6).Value = "right forward"
$ws.Cells.Item(714, 8).Value = 7.28
$ws.Cells.Item(714, 9).Value = 1.49
$ws.Cells.Item(714, 10).Value = 5.76
$ws.Cells.Item(714, 11).Value = 0.92
$ws.Cells.Item(714, 12).Value = 0.41
$ws.Cells.Item(714, 13).Value = 0.15
$ws.Cells.Item(714, 14).Value = 0.04
$ws.Cells.Item(714, 15).Value = 13.0
$ws.Cells.Item(714, 16).Value = 6.61
$ws.Cells.Item(714, 17).Value = 32.75
$ws.Cells.Item(714, 18).Value = 4.53
$ws.Cells.Item(714, 19).Value = 41.0
$ws.Cells.Item(714, 20).Value = 13.0
$ws.Cells.Item(714, 21).Value = 26.0
$ws.Cells.Item(714, 22).Value = 16.0

$ws.Cells.Item(715, 2).Value = 45948.0
$ws.Cells.Item(715, 3).Value = "Global"
$ws.Cells.Item(715, 4).Value = "M"
$ws.Cells.Item(715, 5).Value = "Amir Etien"
$ws.Cells.Item(715, 6).Value = "right forward"
$ws.Cells.Item(715, 8).Value = 8.5
$ws.Cells.Item(715, 9).Value = 1.57
$ws.Cells.Item(715, 10).Value = 6.91
$ws.Cells.Item(715, 11).Value = 0.95
$ws.Cells.Item(715, 12).Value = 0.44
$ws.Cells.Item(715, 13).Value = 0.16
$ws.Cells.Item(715, 14).Value = 0.04
$ws.Cells.Item(715, 15).Value = 13.0
$ws.Cells.Item(715, 16).Value = 5.33
$ws.Cells.Item(715, 17).Value = 33.46
$ws.Cells.Item(715, 18).Value = 4.91
$ws.Cells.Item(715, 19).Value = 40.0
$ws.Cells.Item(715, 20).Value = 11.0
$ws.Cells.Item(715, 21).Value = 35.0
$ws.Cells.Item(715, 22).Value = 14.0

$ws.Cells.Item(716, 2).Value = 45948.0
$ws.Cells.Item(716, 3).Value = "Global"
$ws.Cells.Item(716, 4).Value = "M"
$ws.Cells.Item(716, 5).Value = "Ilan Ihaddadene"
$ws.Cells.Item(716, 6).Value = "center midfield"
$ws.Cells.Item(716, 8).Value = 12.39
$ws.Cells.Item(716, 9).Value = 2.49
$ws.Cells.Item(716, 10).Value = 9.87
$ws.Cells.Item(716, 11).Value = 1.87
$ws.Cells.Item(716, 12).Value = 0.53
$ws.Cells.Item(716, 13).Value = 0.11
$ws.Cells.Item(716, 14).Value = 0.0
$ws.Cells.Item(716, 15).Value = 11.0
$ws.Cells.Item(716, 16).Value = 7.26
$ws.Cells.Item(716, 17).Value = 29.24
$ws.Cells.Item(716, 18).Value = 4.97
$ws.Cells.Item(716, 19).Value = 65.0
$ws.Cells.Item(716, 20).Value = 16.0
$ws.Cells.Item(716, 21).Value = 35.0
$ws.Cells.Item(716, 22).Value = 11.0

$ws.Cells.Item(717, 2).Value = 45948.0
$ws.Cells.Item(717, 3).Value = "Global"
$ws.Cells.Item(717, 4).Value = "M"
$ws.Cells.Item(717, 5).Value = "Emmanuel Valey"
$ws.Cells.Item(717, 6).Value = "left forward"
$ws.Cells.Item(717, 8).Value = 4.42
$ws.Cells.Item(717, 9).Value = 1.07
$ws.Cells.Item(717, 10).Value = 3.34
$ws.Cells.Item(717, 11).Value = 0.71
$ws.Cells.Item(717, 12).Value = 0.26
$ws.Cells.Item(717, 13).Value = 0.1
$ws.Cells.Item(717, 14).Value = 0.01
$ws.Cells.Item(717, 15).Value = 7.0
$ws.Cells.Item(717, 16).Value = 7.44
$ws.Cells.Item(717, 17).Value = 30.62
$ws.Cells.Item(717, 18).Value = 4.21
$ws.Cells.Item(717, 19).Value = 26.0
$ws.Cells.Item(717, 20).Value = 3.0
$ws.Cells.Item(717, 21).Value = 22.0
$ws.Cells.Item(717, 22).Value = 5.0

# Step 5: re-apply the number/alignment formatting used by the rest of the
# table: date format on column B, centred style on column D.
$ws.Cells.Item(697, 2).Copy() | Out-Null
$ws.Range($ws.Cells.Item(698, 2), $ws.Cells.Item(717, 2)).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(697, 4).Copy() | Out-Null
$ws.Range($ws.Cells.Item(698, 4), $ws.Cells.Item(717, 4)).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Step 6: move the selection to match where the editor left off.
$ws.Range("D724").Select() | Out-Null
